$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = 4 (Price), Column E = 5 (Volume(1h))
$D = 4
$E = 5

# Force column D cells to be treated as text so values such as "302.50"
# or "22.30" are not silently converted to numbers and lose their
# trailing zero / formatting (mirrors the original inlineStr text cells).
$priceRows = @(2,3,5,6,9,10,11,14,15,16,17,18,19,20,21,22,23,27,29,30,33,34,35,36,39,40,42,43,44,46,47,48,49,51)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, $D).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Cells.Item(2, $D).Value = "43.057.11"
$ws.Cells.Item(2, $E).Value = "  +1.00%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, $D).Value = "2.351.11"
$ws.Cells.Item(3, $E).Value = "  +2.42%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, $E).Value = "  +0.01%  "

# Row 5 - BNB
$ws.Cells.Item(5, $D).Value = "302.50"
$ws.Cells.Item(5, $E).Value = "  +0.46%  "

# Row 6 - Solana
$ws.Cells.Item(6, $D).Value = "95.40"
$ws.Cells.Item(6, $E).Value = "  -0.28%  "

# Row 7 - XRP
$ws.Cells.Item(7, $E).Value = "  -0.42%  "

# Row 8 - USDC
$ws.Cells.Item(8, $E).Value = "  -0.04%  "

# Row 9 - Cardano
$ws.Cells.Item(9, $D).Value = "0.496"
$ws.Cells.Item(9, $E).Value = "  +0.66%  "

# Row 10 - Avalanche
$ws.Cells.Item(10, $D).Value = "34.02"
$ws.Cells.Item(10, $E).Value = "  -1.57%  "

# Row 11 - Dogecoin
$ws.Cells.Item(11, $D).Value = "0.0787"
$ws.Cells.Item(11, $E).Value = "  +0.16%  "

# Row 12 - Chainlink
$ws.Cells.Item(12, $E).Value = "  -2.57%  "

# Row 13 - TRON
$ws.Cells.Item(13, $E).Value = "  +3.26%  "

# Row 14 - Polkadot
$ws.Cells.Item(14, $D).Value = "6.71"
$ws.Cells.Item(14, $E).Value = "  -0.70%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Cells.Item(15, $D).Value = "2.721.93"
$ws.Cells.Item(15, $E).Value = "  +2.61%  "

# Row 16 - WrappedEther
$ws.Cells.Item(16, $D).Value = "2.361.22"
$ws.Cells.Item(16, $E).Value = "  +3.03%  "

# Row 17 - Polygon
$ws.Cells.Item(17, $D).Value = "0.794"
$ws.Cells.Item(17, $E).Value = "  +1.42%  "

# Row 18 - WrappedBTC
$ws.Cells.Item(18, $D).Value = "43.054.55"
$ws.Cells.Item(18, $E).Value = "  +1.18%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Cells.Item(19, $D).Value = "12.18"
$ws.Cells.Item(19, $E).Value = "  -0.98%  "

# Row 20 - Uniswap
$ws.Cells.Item(20, $D).Value = "6.26"
$ws.Cells.Item(20, $E).Value = "  +4.41%  "

# Row 21 - ShibaInu
$ws.Cells.Item(21, $D).Value = "0.0₃0887"
$ws.Cells.Item(21, $E).Value = "  -0.36%  "

# Row 22 - Litecoin
$ws.Cells.Item(22, $D).Value = "68.05"
$ws.Cells.Item(22, $E).Value = "  +0.43%  "

# Row 23 - BitcoinCash
$ws.Cells.Item(23, $D).Value = "235.09"
$ws.Cells.Item(23, $E).Value = "  +0.06%  "

# Row 24 - ImmutableX
$ws.Cells.Item(24, $E).Value = "  -1.70%  "

# Row 25 - Dai
$ws.Cells.Item(25, $E).Value = "  -0.06%  "

# Row 26 - PancakeSwap
$ws.Cells.Item(26, $E).Value = "  +1.57%  "

# Row 27 - EthereumClassic
$ws.Cells.Item(27, $D).Value = "24.47"
$ws.Cells.Item(27, $E).Value = "  -0.13%  "

# Row 28 - Toncoin
$ws.Cells.Item(28, $E).Value = "  -0.13%  "

# Row 29 - Cosmos
$ws.Cells.Item(29, $D).Value = "9.11"
$ws.Cells.Item(29, $E).Value = "  +0.66%  "

# Row 30 - InjectiveProtocol
$ws.Cells.Item(30, $D).Value = "31.28"
$ws.Cells.Item(30, $E).Value = "  -2.68%  "

# Row 31 - FirstDigitalUSD
$ws.Cells.Item(31, $E).Value = "  +0.01%  "

# Row 32 - Filecoin
$ws.Cells.Item(32, $E).Value = "  +1.03%  "

# Row 33 - Hedera
$ws.Cells.Item(33, $D).Value = "0.0723"
$ws.Cells.Item(33, $E).Value = "  +2.98%  "

# Row 34 - Celestia
$ws.Cells.Item(34, $D).Value = "17.18"
$ws.Cells.Item(34, $E).Value = "  -1.78%  "

# Row 35 - now RenderToken (was ARBITRUM)
$ws.Cells.Item(35, 2).Value = "RenderToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(35, $D).Value = "4.38"
$ws.Cells.Item(35, $E).Value = "  -1.72%  "

# Row 36 - now ARBITRUM (was RenderToken)
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, $D).Value = "1.83"
$ws.Cells.Item(36, $E).Value = "  +4.63%  "

# Row 37 - WEMIXToken
$ws.Cells.Item(37, $E).Value = "  -1.07%  "

# Row 38 - Kaspa
$ws.Cells.Item(38, $E).Value = "  +1.07%  "

# Row 39 - EnergySwap
$ws.Cells.Item(39, $D).Value = "22.30"
$ws.Cells.Item(39, $E).Value = "  +9.51%  "

# Row 40 - LidoDAOToken
$ws.Cells.Item(40, $D).Value = "2.74"
$ws.Cells.Item(40, $E).Value = "  +1.54%  "

# Row 41 - Stellar
$ws.Cells.Item(41, $E).Value = "  -0.23%  "

# Row 42 - Monero
$ws.Cells.Item(42, $D).Value = "103.65"
$ws.Cells.Item(42, $E).Value = "  -37.06%  "

# Row 43 - Maker
$ws.Cells.Item(43, $D).Value = "1.941.72"
$ws.Cells.Item(43, $E).Value = "  -1.14%  "

# Row 44 - VeChain
$ws.Cells.Item(44, $D).Value = "0.0278"
$ws.Cells.Item(44, $E).Value = "  -0.27%  "

# Row 45 - ApeXProtocol
$ws.Cells.Item(45, $E).Value = "  +4.32%  "

# Row 46 - FraxShare
$ws.Cells.Item(46, $D).Value = "9.45"
$ws.Cells.Item(46, $E).Value = "  -9.99%  "

# Row 47 - NEARProtocol
$ws.Cells.Item(47, $D).Value = "2.73"
$ws.Cells.Item(47, $E).Value = "  -1.03%  "

# Row 48 - RocketPoolETH
$ws.Cells.Item(48, $D).Value = "2.584.90"
$ws.Cells.Item(48, $E).Value = "  +2.55%  "

# Row 49 - MultiversX
$ws.Cells.Item(49, $D).Value = "52.80"
$ws.Cells.Item(49, $E).Value = "  -0.70%  "

# Row 50 - HuobiToken
$ws.Cells.Item(50, $E).Value = "  -3.74%  "

# Row 51 - BitcoinSV
$ws.Cells.Item(51, $D).Value = "72.10"
$ws.Cells.Item(51, $E).Value = "  +1.17%  "
